# ---------------------------------------------------------------------------
# Rebuild the "summary" workbook: rename Sheet1 -> block0_trial1, add a new
# sheet block1_trial2, and populate both with a header row (Box_Num,
# Probability_Estimates, Reaction_time, Decision) styled bold/bordered/
# centered, followed by per-trial rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- sheets ----------------------------------------------------------------
$ws1 = $wb.ActiveSheet
$ws1.Name = "block0_trial1"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "block1_trial2"

$headers = @("Box_Num", "Probability_Estimates", "Reaction_time", "Decision")

function Write-BlockSheet {
    param(
        $ws,
        $rows
    )

    # Header row: B1:E1
    $headerRange = $ws.Range("B1:E1")
    $headerRange.Value = $headers

    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108  # xlCenter
    $headerRange.VerticalAlignment = -4160    # xlTop
    $headerRange.Borders.LineStyle = 1        # xlContinuous
    $headerRange.Borders.Weight = 2           # xlThin

    $r = 2
    foreach ($row in $rows) {
        $idxCell = $ws.Cells.Item($r, 1)
        $idxCell.Value = $row[0]
        $idxCell.Font.Bold = $true
        $idxCell.HorizontalAlignment = -4108
        $idxCell.VerticalAlignment = -4160
        $idxCell.Borders.LineStyle = 1
        $idxCell.Borders.Weight = 2

        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $r++
    }
}

$sheet1Rows = @(
    , @(0, 1, 4, 0.8227142999967327, -1)
    , @(1, 2, 3, 0.8621263999993971, "#0000FF")
)

$sheet2Rows = @(
    , @(0, 1, 4, 0.7671456999996735, -1)
    , @(1, 2, 5, 0.6203459999996994, -1)
    , @(2, 3, 5, 1.287835400002223, -1)
    , @(3, 4, 4, 1.044512500000565, -1)
    , @(4, 5, 5, 0.5653464999995776, -1)
    , @(5, 6, 5, 0.8708101999982318, "#0000FF")
)

Write-BlockSheet $ws1 $sheet1Rows
Write-BlockSheet $ws2 $sheet2Rows

$ws1.Activate()
